$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: date "26 February 2021" -> "1 March 2021"
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("26 February 2021", $false, $false, $false, $false, $false, $true, 1, $false, "1 March 2021", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: remove the leading tab before "Purpose. Progress. Brief review..."
# ---------------------------------------------------------------------------
$pPurpose = $d.Paragraphs.Item(7)
$rPurpose = $pPurpose.Range
$tabRange = $d.Range($rPurpose.Start, $rPurpose.Start + 1)
if ($tabRange.Text -eq "`t") {
    $tabRange.Delete()
}

# ---------------------------------------------------------------------------
# Change 3: rewrite the body of the Introduction paragraph (keep leading tab)
# ---------------------------------------------------------------------------
$pIntro = $d.Paragraphs.Item(8)
$rIntro = $pIntro.Range
$introStart = $rIntro.Start
$introLen = $rIntro.End - $rIntro.Start
$introBody = $d.Range($introStart + 1, $introStart + $introLen - 1)

$newIntro = "The purpose of this lab was to write a C program that will be an instruction-level simulator for a limited subset of the ARM instruction set. This program will allow users to tun ARM programs and see the outputs. The objectives of this lab are to introduce the software and process in running code, compiling in C, and introduces the ARM ISA. Several input files are given to test the code that is written. The shell and the simulation routine are the two main sections of the simulator. The goal is to write code for and implement the simulation routine, as the shell code is already given. Code was written within the sim.c and isa.h. Code was written for the data processes, branching, and memory. The results of this code were… All in all, this lab taught students how to design an ARM Architecture Simulator based in the C language."

$introBody.Text = $newIntro

# ---------------------------------------------------------------------------
# Change 4: append new text after the lone tab in the paragraph that follows
# "...How it works. Be specific. <1pg" (Baseline Design section).
# ---------------------------------------------------------------------------
$found = $false
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $core = $p.Range.Text.TrimEnd("`r")
    if ($core -eq "`t") {
        $prev = $d.Paragraphs.Item($i - 1).Range.Text
        if ($prev -like "*Be specific. <1pg*") {
            $found = $true
            $pr = $p.Range
            $tail = $d.Range($pr.End - 1, $pr.End - 1)
            $tail.InsertBefore("The baseline design for this lab was creating code for sim.c and isa.h ")
            break
        }
    }
}
if (-not $found) {
    throw "could not locate the baseline-design tab paragraph"
}
